$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before row 27; existing rows 27-29 shift down to 28-30
$ws.Rows.Item(27).Insert()

# Fill in the new row 27 with data (copy context from surrounding rows, new values per diff)
$ws.Cells.Item(27, 1).Value = 7
$ws.Cells.Item(27, 2).Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(27, 3).Value = "Ñuble"
$ws.Cells.Item(27, 4).Value = 44858
$ws.Cells.Item(27, 4).NumberFormat = $ws.Cells.Item(28, 4).NumberFormat
$ws.Cells.Item(27, 5).Value = 16
$ws.Cells.Item(27, 6).Value = 300000000
$ws.Cells.Item(27, 7).Value = "Espárragos"
$ws.Cells.Item(27, 8).Value = "Sin especificar"
$ws.Cells.Item(27, 9).Value = "Primera"
$ws.Cells.Item(27, 10).Value = 1000
$ws.Cells.Item(27, 11).Value = 1000
$ws.Cells.Item(27, 12).Value = 1200
$ws.Cells.Item(27, 13).Value = 1100
$ws.Cells.Item(27, 14).Value = "$/kilo"
$ws.Cells.Item(27, 15).Value = "Provincia de Diguillín"
$ws.Cells.Item(27, 16).Value = 1100
$ws.Cells.Item(27, 17).Value = 1
$ws.Cells.Item(27, 18).Value = "Hortaliza"

$wb.Save()
